$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: rownames(temp) label + copy of row 1 values, shifted one column right, plus new J value
$ws.Range("A5").Value = "rownames(temp)"
$ws.Range("B5").Value = -282.6592244418332
$ws.Range("C5").Value = 106.95441538661449
$ws.Range("D5").Value = -2.6428008925118998
$ws.Range("E5").Value = 0.0089195273423114934
$ws.Range("F5").Value = -493.65151883511487
$ws.Range("G5").Value = -71.666930048551507
$ws.Range("H5").Value = 187
$ws.Range("I5").Value = 1.972731033408909
$ws.Range("J5").Value = 0

# Row 6: rownames(temp) label + copy of row 2 values, shifted one column right, plus new J value
$ws.Range("A6").Value = "rownames(temp)"
$ws.Range("B6").Value = 3054.9565217391305
$ws.Range("C6").Value = 66.924275134525814
$ws.Range("D6").Value = 45.647958317042686
$ws.Range("E6").Value = [double]"2.5090064061084927e-103"
$ws.Range("F6").Value = 2922.9329272928553
$ws.Range("G6").Value = 3186.9801161854057
$ws.Range("H6").Value = 187
$ws.Range("I6").Value = 1.972731033408909
$ws.Range("J6").Value = 0
